$wb = $excel.ActiveWorkbook

# --- String (play log) cells: append Week 17 data ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "2 4 2 2 -4 2 0 9 2 1 2 25 -2 5 18 3 8 -1 2 -1 14 -2 1 1 14 -3 2 3 1 4 0 3 4 2 4 -3 5 1 1 3 20 0 2 2 1 2 1 10 0 2 4 5 4 8 -2 2 3 11 6 1 3 4 1 6 3 3 3 5 7 6 5 2 4 5 20 -3 3 5 3 0 0 1 3 13 -2 5 18 3 8 -1 -1 4 9 2 1 4 4 1 5 0 -2 1 6 1 25 2 4 1 4 4 5 0 1 6 -1 10 11 2 1 4 9 3 10 5 4 -3 3 3 9 1 2 5 17 16 -3 11 -3 16 6 6 8 5 3 2 1 8 3 5 0 5 2 0 2 2 -1 1 2 1 3 12 10 2 2 7 2 0 5 8 4 5 4 -2 11 0 11 2 0 0 3 1 2 1 13 -1 3 1 5 3 3 11 2 5 5 0 2 9 3 12 26 1 7 4 11 1 0 6 8 4 2 3 3 1 3 2 4 0 1 4 2 4 3 6 4 3 5 2 3 0 12 4 0 3 -2 1 3 1 2 2 4 5 4 3 5 0 4 3 -1 13 3 2 7 0 5 -1 5 4 2 3 6 2 2 0 1 8 3 2 0 8 13 3 5 8 -2 2 3 -2 7 2 -1 3 -2 9 3 4 4 23 8 3 11 13 1 1 2 3 6 -2 2 1 10 3 0 7 2 7 0 0 1 2 3 -3 6 -5 -2 0 2 1 21 5 5 8 6 6 6 12 4 6 3 6 6 5 8 4 2 11 17 2 4 3 8 6 5 8 1 7 0 8 5 14 2 30 -1 5 3 4 1 10 5 3 5 2 3 4 4 4 37"
$ws.Range("B3").Value = "14 6 6 12 9 10 6 19 22 5 9 24 4 -2 5 24 14 1 6 5 2 2 9 7 9 17 6 13 5 41 5 14 3 8 16 52 25 4 5 9 5 4 12 0 11 5 3 0 6 3 7 3 16 12 6 6 3 18 4 6 15 21 16 4 24 0 6 8 3 5 8 7 5 8 11 17 5 29 9 3 8 8 3 9 2 45 6 8 10 4 17 4 6 12 11 11 5 0 4 11 10 7 9 30 3 4 23 50 4 5 13 6 59 5 14 14 3 18 1 22 16 7 9 2 7 9 12 6 4 12 5 6 5 12 9 6 5 2 1 15 14 23 6 12 4 12 7 12 2 3 11 6 12 17 9 3 3 4 9 2 12 7 16 22 4 13 20 24 14 6 50 2 26 9 11 4 0 7 3 1 3 42 10 9 4 10 2 11 4 12 22 13 6 11 4 9 9 -3 7 17 -2 1 7 20 -4 14 6 3 6 5 14 10 9 6 39 10 6 8 9 7 10 8 4 16 0 8 37 22 7 9 10 2 2 6 28 5 7 1 5 16 4 5 32 18 4 5 5 7 11 41 0 11 7 11 2 15 3 10 12 2 30 9 16 13 9 6 8 15 6 5 10 11 3 18 7 12 4 17 5 14 7 7 4 9 40 29 25 9 2 5 9 5 5 7 6 15 6 5 6 6 17 29 5 -2 8 3 14 1 37 30 6 15 6 38 9 1 13 3 10 12 16 4 6 3 -6 11 14 7 19 5 5 17 12 7 5 8 4 18 9 5 -3 24 10 0 1 3 5 10 6 2 17 4 8 2 5 6 6 3 3 7 4 6 12 6 -4 8 -2 6 0 5 11 7 3 5 0 13 1 12 3 9 1 8 4 3 10 12 2"
$ws.Range("C2").Value = "8 2 6 9 1 8 -3 1 3 1 3 3 3 7 3 2 2 11 2 6 15 25 2 3 2 4 1 0 4 5 -2 2 1 -2 4 -2 1 1 -4 11 5 2 -4 3 13 1 3 1 2 6 4 -2 6 27 1 -1 3 11 11 4 3 2 8 1 1 5 2 7 -1 4 4 3 0 7 3 6 4 4 -2 15 9 6 5 3 6 3 1 -3 3 25 1 -1 5 4 5 6 2 -1 3 3 1 5 2 -4 8 0 49 2 7 1 1 2 3 10 -1 11 4 5 7 7 6 -1 5 11 21 5 7 7 5 0 2 0 3 1 14 3 26 5 10 -1 -1 6 1 -4 14 4 2 5 1 3 1 0 11 21 2 4 0 10 -2 2 0 5 2 12 1 -1 2 3 8 21 2 8 1 6 2 8 1 0 6 2 2 3 2 10 3 1 12 5 2 15 3 8 1 7 9 3 0 0 4 28 15 4 2 5 21 -1 6 16 1 12 14 42 4 1 4 1 0 -5 0 4 -4 2 6 5 3 0 1 -2 8 4 2 3 7 6 10 4 0 18 2 18 2 12 6 -2 13 4 5 2 5 36 5 25 3 7 2 3 6 3 8 3 1 6 4 5 9 8 9 0 9 5 1 8 9 4 -1 0 1 32 1 1 2 3 13 2 3 -2 -1 4 1 12 5 18 -1 3 5 3 8 -2 4 5 3 8 6 8 0 -1 -2 13 3 2 2 3 1 20 4 30 14 3 3 16 15 3 29 9 3 0 1 7 7 9 8 5 -4 7 14 -3 5 -2 8 5 4 1 11 2 4 -1 0 2 3 6 2 1 3 8 9 2 1 2 3 6 0 3 9 1 2 0 9 12 1 6 3 10 7 2 4 4 15 4 5 2 11 0 20 6 7 3 1 4 3 0 3 4 0 1 3 8 12 5 5 2 3 2 7 2 6 3 8 3 10 4 6 7 2 7 4 0 1 0 0 1 5 -2 3 -2 32 2 6 4 1 2 4 1 3 4 7 3 6 0 -2 4 2 14 2"
$ws.Range("C3").Value = "7 -2 6 5 10 16 13 4 7 12 10 5 37 10 8 3 7 9 9 12 7 4 6 12 -1 13 10 13 5 9 3 11 10 17 13 10 25 6 4 3 8 9 14 18 9 23 30 9 13 2 13 16 7 4 10 61 25 12 15 7 17 7 19 34 14 8 14 5 9 4 7 5 21 7 16 4 8 12 6 6 23 7 19 10 25 26 11 1 28 2 11 5 10 9 6 14 12 12 3 12 19 10 9 15 9 39 23 6 9 5 23 6 -1 2 27 12 -3 12 14 0 9 14 41 1 3 -4 7 7 7 2 11 7 12 21 10 6 16 12 8 14 6 19 8 7 11 15 21 8 11 24 9 6 11 9 4 3 6 25 8 21 17 50 10 11 20 28 5 9 39 16 11 12 2 -3 7 5 19 0 5 5 9 4 1 30 8 2 22 5 11 6 17 5 7 10 18 10 14 9 10 18 6 30 14 17 8 5 7 17 3 11 12 13 5 17 53 5 9 5 16 4 12 32 11 16 7 15 16 6 0 -4 0 -7 16 12 19 14 16 18 0 9 9 8 2 10 29 3 10 15 15 3 17 6 5 22 18 13 5 6 8 5 8 14 18 22 19 17 13 0 1 62 12 17 7 5 11 12 1 4 8 5 18 7 9 13 11 5 -4 -3 3 5 12 5 9 8 10 11 11 9 13 -1 5 18 18 5 7 3 14 8 50 3 8 9 4 6 9 16 22 20 9 22 7 7 39 3 21 8 8 1 8 9 1 7"

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B4").Value = "64 63 58 55 60 68 66 64 65 50 65 58 68 63 63 46 55 68 70 59 73 63 57 67 63 59 60 67 61 62 67 66 62 62 64 59 60 62 65 59 60 68 66"
$ws.Range("B5").Value = "75 23 19 26 22 23 30 21 16 10 22 18 34 38 25 0 10 25 16 19 26 32 17 18 29 38 28 28 23 22 12 26 35 0 21 19 20 21 21 29 21 27 21"
$ws.Range("D3").Value = "41 51 39 30 46 38 59 42 58 45 50 27 40 45 20 37 56 63 45 38 58 56 32 50 56 38 56 44 38 38 32 46 41 39 39 61 46 39 45 50 46 31 37 42 38 42 38 64 43 47 25 40 49 45 49 41 34 27 51 32 40 49 58 63 53 43 44 53 51"
$ws.Range("D4").Value = "0 1 0 0 0 0 10 0 16 16 -7 0 0 9 0 0 13 13 0 3 14 0 0 3 0 0 1 0 0 0 0 0 0 0 0 48 5 0 0 18 -1 0 0 17 0 0 0 0 12 0 0 16 12 0 0 15 0 0 55 0 0 0 18 0 4 0 9 9 7"
$ws.Range("D5").Value = "0 5 0 15 0 0 15 14 0 0 0 0 0 0 14 0 0 0 0 6 0 0 0 0 0 14 9 0 10 19 0 14 6 3 0 13 13 0 0 5 13 0 10 0 1 14 0 3 0 14 5 6 8 0 8 0 17 11 11 0 20 0 0 15 10 5 0 14"

# --- Numeric summary cells ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 188
$ws.Range("F2").Value = 61
$ws.Range("G2").Value = 45
$ws.Range("J2").Value = 24
$ws.Range("L2").Value = 348
$ws.Range("M2").Value = 221
$ws.Range("O2").Value = 19
$ws.Range("P2").Value = 9
$ws.Range("Q2").Value = 599
$ws.Range("C3").Value = 188
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 37
$ws.Range("F3").Value = 124
$ws.Range("G3").Value = 51
$ws.Range("H3").Value = 33
$ws.Range("I3").Value = 74
$ws.Range("J3").Value = 76
$ws.Range("N3").Value = 18
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 215
$ws.Range("F2").Value = 65
$ws.Range("G2").Value = 59
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 31
$ws.Range("L2").Value = 285
$ws.Range("M2").Value = 181
$ws.Range("O2").Value = 22
$ws.Range("P2").Value = 12
$ws.Range("Q2").Value = 584
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 166
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 36
$ws.Range("F3").Value = 98
$ws.Range("G3").Value = 35
$ws.Range("H3").Value = 41
$ws.Range("I3").Value = 61
$ws.Range("J3").Value = 45
$ws.Range("N3").Value = 37
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 83
$ws.Range("D2").Value = 69
$ws.Range("F2").Value = 244
$ws.Range("G2").Value = 230
$ws.Range("J2").Value = 108
$ws.Range("K2").Value = 102
$ws.Range("L2").Value = 74
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 21
$ws.Range("O2").Value = 17
$ws.Range("B3").Value = 40
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 7
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 10
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 20
$ws.Range("D4").Value = 5
